$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet
$ws.Name = "25preds_et"

# Update HT (half-time) model prediction values in columns B, C, D
$ws.Range("C2").Value = 7.5391666666666604
$ws.Range("D2").Value = 6.2850000000000001
$ws.Range("C3").Value = 12.3808333333333
$ws.Range("D3").Value = 12.04
$ws.Range("C4").Value = 12.37
$ws.Range("D4").Value = 10.6625
$ws.Range("C5").Value = 18.074166666666599
$ws.Range("D5").Value = 16.199166666666599
$ws.Range("C6").Value = 15.4308333333333
$ws.Range("D6").Value = 12.011666666666599
$ws.Range("B7").Value = 29.520833339617798
$ws.Range("C7").Value = 2.2799999999999998
$ws.Range("D7").Value = 14.68
$ws.Range("C8").Value = 11.8716666666666
$ws.Range("D8").Value = 12.22
$ws.Range("C9").Value = 13.9166666666666
$ws.Range("D9").Value = 6.0333333333333297
$ws.Range("C11").Value = 7.55416666666666
$ws.Range("D11").Value = 6.5233333333333299
$ws.Range("C12").Value = 11.803333333333301
$ws.Range("D12").Value = 15.365
$ws.Range("C14").Value = 10.5283333333333
$ws.Range("D14").Value = 13.1183333333333
$ws.Range("C15").Value = 18.0066666666666
$ws.Range("D15").Value = 14.209166666666601
$ws.Range("C16").Value = 8.5733333333333306
$ws.Range("D16").Value = 9.7491666666666603
$ws.Range("C17").Value = 14.289166666666601
$ws.Range("D17").Value = 9.6724999999999994
$ws.Range("C18").Value = 10.7116666666666
$ws.Range("D18").Value = 13.897500000000001
$ws.Range("C19").Value = 1.98166666666666
$ws.Range("D19").Value = 14.5541666666666
$ws.Range("C20").Value = 16.3258333333333
$ws.Range("D20").Value = 10.376666666666599
$ws.Range("B21").Value = 17.1945913898806
$ws.Range("C21").Value = 12.4825
$ws.Range("D21").Value = 8.4866666666666593
$ws.Range("C22").Value = 13.3433333333333
$ws.Range("D22").Value = 7.0341666666666596
$ws.Range("C23").Value = 20.816666666666599
$ws.Range("D23").Value = 15.3125
$ws.Range("C24").Value = 8.3149999999999995
$ws.Range("D24").Value = 8.7541666666666593
$ws.Range("C25").Value = 11.9541666666666
$ws.Range("D25").Value = 9.9891666666666605
$ws.Range("C26").Value = 15.7183333333333
$ws.Range("D26").Value = 4.3841666666666601
$ws.Range("C27").Value = 7.7066666666666599
$ws.Range("D27").Value = 8.7633333333333301
$ws.Range("C28").Value = 8.0933333333333302
$ws.Range("D28").Value = 14.4433333333333
$ws.Range("C29").Value = 11.6416666666666
$ws.Range("D29").Value = 8.8783333333333303
$ws.Range("C30").Value = 12.1241666666666
$ws.Range("D30").Value = 17.273333333333301
$ws.Range("C31").Value = 14.220833333333299
$ws.Range("D31").Value = 16.739999999999998
$ws.Range("C32").Value = 10.101666666666601
$ws.Range("D32").Value = 12.1516666666666
$ws.Range("C33").Value = 9.5
$ws.Range("D33").Value = 14.758333333333301
$ws.Range("C34").Value = 8.6983333333333306
$ws.Range("D34").Value = 7.4858333333333302
$ws.Range("C35").Value = 13.5591666666666
$ws.Range("D35").Value = 15.584166666666601
$ws.Range("C36").Value = 14.5691666666666
$ws.Range("D36").Value = 14.1833333333333
$ws.Range("A37").Value = 36.570265762019197
$ws.Range("C37").Value = 20.884166666666601
$ws.Range("D37").Value = 15.6225
$ws.Range("C38").Value = 13.8116666666666
$ws.Range("D38").Value = 9.2916666666666607
$ws.Range("C39").Value = 9.8066666666666595
$ws.Range("D39").Value = 12.101666666666601
$ws.Range("C40").Value = 10.348333333333301
$ws.Range("D40").Value = 12.0891666666666
$ws.Range("C41").Value = 12.2383333333333
$ws.Range("D41").Value = 11.011666666666599
